# Insert a new weekly observation row for "Terminal La Palmera de La Serena - Zanahoria"
# at row 601, pushing the existing rows 601:645 down to 602:646.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("601").Insert()

$ws.Range("A601").Value = 8
$ws.Range("B601").Value = "Terminal La Palmera de La Serena"
$ws.Range("C601").Value = "Coquimbo"
$ws.Range("D601").Value = 45265
$ws.Range("E601").Value = 4
$ws.Range("F601").Value = 100114013
$ws.Range("G601").Value = "Zanahoria"
$ws.Range("H601").Value = "Sin especificar"
$ws.Range("I601").Value = "Primera"
$ws.Range("J601").Value = 460
$ws.Range("K601").Value = 5500
$ws.Range("L601").Value = 6000
$ws.Range("M601").Value = 5750
$ws.Range("N601").Value = "`$/saco 20 kilos"
$ws.Range("O601").Value = "Provincia del Elquí"
$ws.Range("P601").Value = 288
$ws.Range("Q601").Value = 20
$ws.Range("R601").Value = "Hortaliza"
